$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.215.57"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "2.917.36"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'371.36"
$ws.Range("E5").Value = "  +3.98%  "
$ws.Range("D6").Value = "'104.61"
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("D7").Value = "'0.542"
$ws.Range("E7").Value = "  -5.43%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  -5.87%  "
$ws.Range("D10").Value = "'37.19"
$ws.Range("E10").Value = "  -4.71%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "'0.0838"
$ws.Range("E12").Value = "  -4.43%  "
$ws.Range("D13").Value = "'18.48"
$ws.Range("E13").Value = "  -5.43%  "
$ws.Range("D14").Value = "3.384.14"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "'7.40"
$ws.Range("E15").Value = "  -4.77%  "
$ws.Range("D16").Value = "2.929.94"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'0.951"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("D18").Value = "51.257.60"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "'3.32"
$ws.Range("E19").Value = "  -5.64%  "
$ws.Range("D20").Value = "'7.29"
$ws.Range("E20").Value = "  -4.65%  "
$ws.Range("D21").Value = "'13.04"
$ws.Range("E21").Value = "  -6.22%  "
$ws.Range("D22").Value = "0.0₃0946"
$ws.Range("E22").Value = "  -3.79%  "
$ws.Range("D23").Value = "'68.43"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").Value = "'260.88"
$ws.Range("E24").Value = "  -3.71%  "
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("E26").Value = "  +4.03%  "
$ws.Range("D27").Value = "'0.174"
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'7.40"
$ws.Range("E29").Value = "  -5.20%  "
$ws.Range("D30").Value = "'25.95"
$ws.Range("E30").Value = "  -4.05%  "
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "'9.94"
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").Value = "'6.13"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").Value = "'35.28"
$ws.Range("E34").Value = "  -6.72%  "
$ws.Range("D35").Value = "'2.13"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").Value = "'50.68"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'0.0424"
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("D39").Value = "'3.11"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("D40").Value = "'2.70"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("D41").Value = "'17.01"
$ws.Range("E41").Value = "  -6.70%  "
$ws.Range("D42").Value = "'1.86"
$ws.Range("E42").Value = "  -6.35%  "
$ws.Range("E43").Value = "  -5.72%  "
$ws.Range("D44").Value = "'22.27"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("D45").Value = "'118.21"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").Value = "'2.10"
$ws.Range("E46").Value = "  -3.69%  "
$ws.Range("D47").Value = "2.052.36"
$ws.Range("E47").Value = "  -4.12%  "
$ws.Range("E48").Value = "  -5.99%  "
$ws.Range("D49").Value = "'3.20"
$ws.Range("E49").Value = "  -7.28%  "
$ws.Range("D50").Value = "3.210.19"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").Value = "'0.239"
$ws.Range("E51").Value = "  -3.63%  "
